$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Case_1")

$ws.Range("B2").Value = 2071.1659551
$ws.Range("D2").Value = 0.004050567304590923

$ws.Range("B3").Value = 2916.6421686
$ws.Range("D3").Value = -0.1055165180861513

$ws.Range("B4").Value = 3574.7105138
$ws.Range("D4").Value = -0.001086947267102302
